# Edit script: reproduce the commit's two changes.
#
# 1) Slide 6's table (graphicFrame) switches its table style (tableStyleId)
#    from {D3A3BA8E-F770-495B-80A6-9400A1A838F2} to
#    {CD1D2B7C-0404-4998-83B7-B4A55C9DCC7E}.
#
# 2) The deck's applied theme colour scheme changes from the "Integral"
#    palette to the default "Office" palette (the underlying OOXML diff
#    shows the theme1.xml/theme2.xml parts swapping contents end to end;
#    the part that is reachable/meaningful through the PowerPoint object
#    model is the active theme's 12 colour-scheme slots, which this
#    recreates by setting each ThemeColorScheme entry to the stock
#    "Office" theme colour values).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl = $s6.Shapes.Item(2).Table
$tbl.ApplyStyle("{CD1D2B7C-0404-4998-83B7-B4A55C9DCC7E}")

# --- 2. Swap the active theme palette from "Integral" to "Office" --------
# ThemeColorScheme indices follow MsoThemeColorSchemeIndex ordering:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB is stored as a BGR-packed integer (0xBBGGRR), as in the real
# PowerPoint object model.
$cs = $p.Slides.Item(1).ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
